$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - SOC_XSIT_1 (already present in A2)
$ws.Range("B2").Value = Get-Date -Year 2014 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2").Value = Get-Date -Year 2009 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = 4.8
$ws.Range("G2").Value = "social"
$ws.Range("H2").Value = 4

# Row 3 - SOC_XSIT_2
$ws.Range("A3").Value = "SOC_XSIT_2"
$ws.Range("B3").Value = Get-Date -Year 2014 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("C3").Value = Get-Date -Year 2011 -Month 3 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = "M"
$ws.Range("F3").Value = 3.3
$ws.Range("G3").Value = "social"
$ws.Range("H3").Value = 3

# Row 4 - SOC_XSIT_3
$ws.Range("A4").Value = "SOC_XSIT_3"
$ws.Range("B4").Value = Get-Date -Year 2014 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = Get-Date -Year 2008 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("E4").Value = "F"
$ws.Range("F4").Value = 5.6
$ws.Range("G4").Value = "social"
$ws.Range("H4").Value = 5

# Row 5 - SOC_XSIT_4
$ws.Range("A5").Value = "SOC_XSIT_4"
$ws.Range("B5").Value = Get-Date -Year 2014 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("C5").Value = Get-Date -Year 2009 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("E5").Value = "F"
$ws.Range("F5").Value = 4.7
$ws.Range("G5").Value = "social"
$ws.Range("H5").Value = 4

# Update active cell selection to C14, matching sheetView change
$ws.Range("C14").Select()
